$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 28990
$ws.Range("J105").Value = 28990
$ws.Range("L105").Value = 28990
$ws.Range("N105").Value = -35978

$ws.Range("H129").Value = 1215.5217
$ws.Range("I129").Value = 677.6667
$ws.Range("J129").Value = 1253.0465
$ws.Range("K129").Value = 2033.0001
$ws.Range("L129").Value = 3759.1395
$ws.Range("M129").Value = 2966.9999
$ws.Range("N129").Value = -13759.1395

$ws.Range("H137").Value = 759154.7
$ws.Range("I137").Value = 2385265.5
$ws.Range("J137").Value = 2824.0698
$ws.Range("K137").Value = 7155796.5
$ws.Range("L137").Value = 8472.2094
$ws.Range("M137").Value = -7153246.5
$ws.Range("N137").Value = -13572.2094

$ws.Range("H138").Value = 3480.17
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 3480.17
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 10440.51
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -20720.51

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4256.5713
$ws.Range("I32").Value = 4568.5317
$ws.Range("K32").Value = 4568.5317
$ws.Range("M32").Value = -4281.5317

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H137").Value = 42097.5
$ws.Range("J137").Value = 42097.5
$ws.Range("L137").Value = 42097.5
$ws.Range("N137").Value = -52297.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 41980
$ws.Range("J130").Value = 41980
$ws.Range("L130").Value = 41980
$ws.Range("N130").Value = -52020

$ws.Range("H132").Value = 53958.355
$ws.Range("J132").Value = 53958.355
$ws.Range("L132").Value = 53958.355
$ws.Range("N132").Value = -64078.355

$ws.Range("H135").Value = 49951
$ws.Range("J135").Value = 49951
$ws.Range("L135").Value = 49951
$ws.Range("N135").Value = -60091

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 41000
$ws.Range("J138").Value = 41000
$ws.Range("L138").Value = 41000
$ws.Range("N138").Value = -51280

$ws.Range("H140").Value = 48316.96
$ws.Range("J140").Value = 48316.96
$ws.Range("L140").Value = 48316.96
$ws.Range("N140").Value = -58676.96

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 194712.86
$ws.Range("I31").Value = 365976.88
$ws.Range("J31").Value = 2689.5757
$ws.Range("K31").Value = 365976.88
$ws.Range("L31").Value = 2689.5757
$ws.Range("M31").Value = -365681.88
$ws.Range("N31").Value = -3279.5757

$ws.Range("H34").Value = 194712.86
$ws.Range("I34").Value = 365976.88
$ws.Range("J34").Value = 2689.5757
$ws.Range("K34").Value = 365976.88
$ws.Range("L34").Value = 2689.5757
$ws.Range("M34").Value = -365774.88
$ws.Range("N34").Value = -3093.5757

$ws.Range("H115").Value = 24949.5
$ws.Range("I115").Value = 24999
$ws.Range("J115").Value = 24900
$ws.Range("K115").Value = 24999
$ws.Range("L115").Value = 24900
$ws.Range("M115").Value = -23824
$ws.Range("N115").Value = -27250

$ws.Range("H132").Value = 4474.9
$ws.Range("I132").Value = 3676.4707
$ws.Range("J132").Value = 8999.333000000001
$ws.Range("K132").Value = 11029.4121
$ws.Range("L132").Value = 26997.999
$ws.Range("M132").Value = -8499.4121
$ws.Range("N132").Value = -32057.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1320.5922
$ws.Range("I68").Value = 920.5172
$ws.Range("J68").Value = 1567.4468
$ws.Range("K68").Value = 2761.5516
$ws.Range("L68").Value = 4702.3404
$ws.Range("M68").Value = -1950.5516
$ws.Range("N68").Value = -6324.3404

$ws.Range("H71").Value = 1320.5922
$ws.Range("I71").Value = 920.5172
$ws.Range("J71").Value = 1567.4468
$ws.Range("K71").Value = 8284.6548
$ws.Range("L71").Value = 14107.0212
$ws.Range("M71").Value = -4228.6548
$ws.Range("N71").Value = -22219.0212

$ws.Range("H121").Value = 1189.4912
$ws.Range("I121").Value = 1696.3334
$ws.Range("J121").Value = 1161.3334
$ws.Range("K121").Value = 5089.0002
$ws.Range("L121").Value = 3484.0002
$ws.Range("M121").Value = -3779.0002
$ws.Range("N121").Value = -6104.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H46").Value = 23887.467
$ws.Range("J46").Value = 24927.846
$ws.Range("L46").Value = 24927.846
$ws.Range("N46").Value = -25239.846

$ws.Range("H113").Value = 2802.75
$ws.Range("I113").Value = 2988.8572
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 2988.8572
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -818.8571999999999
$ws.Range("N113").Value = -5840

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H122").Value = 4062.3845
$ws.Range("I122").Value = 2781.3
$ws.Range("J122").Value = 8332.666999999999
$ws.Range("K122").Value = 8343.900000000001
$ws.Range("L122").Value = 24998.001
$ws.Range("M122").Value = -5893.900000000001
$ws.Range("N122").Value = -29898.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 42000
$ws.Range("J47").Value = 42000
$ws.Range("L47").Value = 42000
$ws.Range("N47").Value = -42980

$ws.Range("H52").Value = 42000
$ws.Range("J52").Value = 42000
$ws.Range("L52").Value = 42000
$ws.Range("N52").Value = -42466

$ws.Range("H115").Value = 35633.332
$ws.Range("J115").Value = 35633.332
$ws.Range("L115").Value = 35633.332
$ws.Range("N115").Value = -37983.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 86193790
$ws.Range("I96").Value = 144357940
$ws.Range("J96").Value = 4763985
$ws.Range("K96").Value = 144357940
$ws.Range("L96").Value = 4763985
$ws.Range("M96").Value = -4766731
$ws.Range("N96").Value = -4766731

Write-Output "Applied all Chocobo_Profits updates"